# RBA v2.4 - Alteracao da pasta de origem
#
# Applies the text replacements described by the diff:
#   - Body: " 000110318713 - 2 " -> "  " (student RA number cleared to blank)
#   - Body: "QWR" -> "TERE" (directoria placeholder in the "A ..." sentence)
#   - Header: "QWER" -> "TRE"
#   - Header: "QWR"  -> "TERE"
#   - Header: "Qwer" -> "Tre" (5x)
#   - Header: "qwer" -> "tre" (3x)

$d = $word.ActiveDocument

function Replace-ExactText($range, [string]$findText, [string]$replaceText) {
    $f = $range.Find
    $f.ClearFormatting()
    $f.Replacement.ClearFormatting()
    $f.Text = $findText
    $f.Replacement.Text = $replaceText
    $f.Forward = $true
    $f.Wrap = 0
    $f.Format = $false
    $f.MatchCase = $true
    $f.MatchWholeWord = $false
    $f.MatchWildcards = $false
    $f.MatchSoundsLike = $false
    $f.MatchAllWordForms = $false
    $null = $f.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

# --- Body (main document story) ---------------------------------------
Replace-ExactText $d.Content " 000110318713 - 2 " "  "
Replace-ExactText $d.Content "QWR" "TERE"

# --- Header (primary header story) -------------------------------------
$headerRange = $null
foreach ($story in $d.StoryRanges) {
    if ($story.StoryType -eq 7) {
        $headerRange = $story
    }
}

Replace-ExactText $headerRange "QWER" "TRE"
Replace-ExactText $headerRange "QWR" "TERE"
Replace-ExactText $headerRange "Qwer" "Tre"
Replace-ExactText $headerRange "qwer" "tre"
